$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refreshed TPM-derived numbers (only the cells that actually changed) ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2063233333333333
$ws.Range("H2").Value = 0.61897
$ws.Range("M2").Value = 0.3987243333333333
$ws.Range("N2").Value = 1.196173
$ws.Range("O2").Value = 0.5027922793035905
$ws.Range("P2").Value = 0.5027922793035905
$ws.Range("Q2").Value = 0.08226613353444444
$ws.Range("R2").Value = 0.74039520181
$ws.Range("S2").Value = 0.5027922793035905
$ws.Range("T2").Value = 0.5027922793035905

# --- Row 3: refreshed TPM-derived numbers (only the cells that actually changed) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2063233333333333
$ws.Range("H3").Value = 0.61897
$ws.Range("M3").Value = 0.321599
$ws.Range("N3").Value = 0.964797
$ws.Range("O3").Value = 0.4055370608559684
$ws.Range("P3").Value = 0.4055370608559684
$ws.Range("Q3").Value = 0.06635337767666667
$ws.Range("R3").Value = 0.59718039909
$ws.Range("S3").Value = 0.4055370608559684
$ws.Range("T3").Value = 0.4055370608559684

# --- New row 4: MuSCs -> Dkk4 -> Kremen2 -> Resolving-Mac ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Dkk4"
$ws.Range("C4").Value = "Kremen2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2063233333333333
$ws.Range("H4").Value = 0.61897
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07269666666666667
$ws.Range("N4").Value = 0.21809
$ws.Range("O4").Value = 0.09167065984044119
$ws.Range("P4").Value = 0.09167065984044119
$ws.Range("Q4").Value = 0.01499901858888889
$ws.Range("R4").Value = 0.1349911673
$ws.Range("S4").Value = 0.09167065984044119
$ws.Range("T4").Value = 0.09167065984044119
